$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.888.78'
$ws.Range("E2").Value = '  -0.28%  '
$ws.Range("D3").Value = '1.629.93'
$ws.Range("E3").Value = '  -0.74%  '
$ws.Range("E4").Value = '  -0.23%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.17'
$ws.Range("E5").Value = '  -0.78%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.522'
$ws.Range("E6").Value = '  -0.48%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.25%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.41'
$ws.Range("E8").Value = '  -0.52%  '
$ws.Range("E9").Value = '  -0.40%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0881'
$ws.Range("E11").Value = '  -0.23%  '
$ws.Range("D12").Value = '1.861.32'
$ws.Range("E12").Value = '  -0.75%  '
$ws.Range("D13").Value = '1.625.43'
$ws.Range("E13").Value = '  -0.81%  '
$ws.Range("E14").Value = '  -1.49%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.562'
$ws.Range("E15").Value = '  -1.73%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.29'
$ws.Range("E16").Value = '  -0.35%  '
$ws.Range("D17").Value = '27.891.77'
$ws.Range("E17").Value = '  -0.27%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '229.28'
$ws.Range("E18").Value = '  -1.57%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.67'
$ws.Range("E19").Value = '  +0.91%  '
$ws.Range("E20").Value = '  -0.34%  '
$ws.Range("E21").Value = '  -0.32%  '
$ws.Range("E22").Value = '  -0.82%  '
$ws.Range("E23").Value = '  -4.11%  '
$ws.Range("E24").Value = '  -1.35%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.06'
$ws.Range("E25").Value = '  +0.66%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.89'
$ws.Range("E26").Value = '  -0.09%  '
$ws.Range("E27").Value = '  -0.18%  '
$ws.Range("E28").Value = '  -1.16%  '
$ws.Range("E29").Value = '  -0.18%  '
$ws.Range("E30").Value = '  -0.88%  '
$ws.Range("E31").Value = '  -0.71%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.42'
$ws.Range("E32").Value = '  +0.60%  '
$ws.Range("E33").Value = '  +0.62%  '
$ws.Range("D34").Value = '1.392.01'
$ws.Range("E34").Value = '  -1.23%  '
$ws.Range("E35").Value = '  +0.13%  '
$ws.Range("E36").Value = '  +9.34%  '
$ws.Range("E37").Value = '  -1.02%  '
$ws.Range("E38").Value = '  +1.04%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.557'
$ws.Range("E39").Value = '  -1.23%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.851'
$ws.Range("E40").Value = '  -3.32%  '
$ws.Range("B41").Value = 'WEMIXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.02'
$ws.Range("E41").Value = '  -1.16%  '
$ws.Range("B42").Value = 'PaxDollar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.998'
$ws.Range("E42").Value = '  -0.26%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.84'
$ws.Range("E43").Value = '  -1.50%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '65.72'
$ws.Range("E44").Value = '  -2.16%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.44'
$ws.Range("E45").Value = '  -1.95%  '
$ws.Range("D46").Value = '1.771.73'
$ws.Range("E47").Value = '  -2.93%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '88.33'
$ws.Range("E48").Value = '  +0.40%  '
$ws.Range("E49").Value = '  +1.38%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0504'
$ws.Range("E50").Value = '  -0.49%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.62'
$ws.Range("E51").Value = '  -0.05%  '
